$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set new shared-string cells first, in the exact order the strings are
# --- first introduced, so the rebuilt sharedStrings table lines up with the
# --- target workbook (unchanged strings keep their relative order, removed
# --- ones are pruned, and brand-new ones are appended in order of first use).
$ws.Range("L1").Value = "ExpectedNumOfTxn"
$ws.Range("M1").Value = "AC Screen list"
$ws.Range("I2").Value = "8-02-200000"
$ws.Range("N1").Value = "RelationCode"
$ws.Range("H2").Value = "8-02"
$ws.Range("O2").Value = "Below 1M"
$ws.Range("O1").Value = "TurnoverM"
$ws.Range("P1").Value = "TurnoverA"
$ws.Range("P2").Value = "1M to 5M"
$ws.Range("Q1").Value = "debitTxnNum"
$ws.Range("R1").Value = "TurnoverDebitMonth"

# --- Remaining header cells (unchanged text, kept for completeness) ---
$ws.Range("A1").Value = "Cid"
$ws.Range("B1").Value = "CategoryProduct"
$ws.Range("C1").Value = "Acc name"
$ws.Range("D1").Value = "Acc name2"
$ws.Range("E1").Value = "SignOffData"
$ws.Range("F1").Value = "sbpCompany"
$ws.Range("G1").Value = "sbpSector"
$ws.Range("H1").Value = "sbpSubSector"
$ws.Range("I1").Value = "sbpSegment"
$ws.Range("J1").Value = "Jholder"
$ws.Range("K1").Value = "Purpose"

# --- Row 2 data ---
$ws.Range("A2").Value = 10001231
$ws.Range("B2").Value = "1-011"
$ws.Range("C2").Value = "abc"
$ws.Range("D2").Value = "abc"

$ws.Range("E2").Value = 20230106
$ws.Range("E2").NumberFormat = "@"

$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 8

$ws.Range("H2").NumberFormat = "@"

$ws.Range("J2").Value = 16206304
$ws.Range("K2").Value = "testing"
$ws.Range("L2").Value = 20
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 6
$ws.Range("Q2").Value = 20
$ws.Range("R2").Value = "Below 1M"

# --- View state ---
$ws.Range("O1").Select() | Out-Null
